$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark "training calendar view" row as delivered (F130: 0 -> 1) ---
$ws.Range("F130").Value = 1

# --- Append two new backlog rows (137, 138) ---
# Copy formatting/styles from an existing visible data row (130) so the
# new rows pick up the same cell styles (s="3" for title, s="4" for the rest)
# without Excel re-mapping them to a different (but equivalent) style index.
$ws.Rows("130:130").Copy()
$ws.Rows("137:137").Insert(-4121)
$ws.Rows("130:130").Copy()
$ws.Rows("138:138").Insert(-4121)

$ws.Range("A137").Value = "افزودن قابليت صفحه بندي به ليست هاي موجود در نرم افزار"
$ws.Range("B137").Value = "دوم"
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 0
$ws.Range("E137").Value = 0
$ws.Range("F137").Value = 0

$ws.Range("A138").Value = "قابليت ورود به نرم افزار و افزودن سطح دسترسي به بخش هاي مختلف"
$ws.Range("B138").Value = "دوم"
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 0
$ws.Range("E138").Value = 0
$ws.Range("F138").Value = 0

# --- Grow Table1 / its AutoFilter to cover the two new rows ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F138"))

# --- Restore the active cell/selection recorded in the saved workbook ---
$null = $ws.Range("F98").Select()
